$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in B2:C6 per the new run
$ws.Range("B2").Value = 0.01396621389009421
$ws.Range("C2").Value = 0.03310379733627047

$ws.Range("B3").Value = 6.853083661559884
$ws.Range("C3").Value = 5.74732892315861

$ws.Range("B4").Value = 14.77416229248047
$ws.Range("C4").Value = 14.35161113739014

# Row 5 values change AND the bold styling moves from B5 to C5
$ws.Range("B5").Value = 0.08938261490520505
$ws.Range("C5").Value = 0.09419863735753496

# Swap formatting: B5 becomes plain (like B2), C5 becomes bold (like C2)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null

$ws.Range("B6").Value = 1.38886579167174
$ws.Range("C6").Value = 1.101778914883925
